$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '30.577.05'
$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  -0.22%  '

# Row 3
$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '1.893.50'
$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +0.39%  '

# Row 4
$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  -0.07%  '

# Row 5
$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '238.13'
$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +1.32%  '

# Row 6
$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  -0.01%  '

# Row 7
$c = $ws.Range('D7')
$c.NumberFormat = '@'
$c.Value = '0.4896'
$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.47%  '

# Row 8
$c = $ws.Range('D8')
$c.NumberFormat = '@'
$c.Value = '0.2937'
$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  +1.81%  '

# Row 9
$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '0.06698'
$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.67%  '

# Row 10
$c = $ws.Range('D10')
$c.NumberFormat = '@'
$c.Value = '1.910.99'
$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +1.34%  '

# Row 11
$c = $ws.Range('D11')
$c.NumberFormat = '@'
$c.Value = '17.19'
$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  +3.03%  '

# Row 12
$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '0.07340'
$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +1.51%  '

# Row 13
$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '5.152'
$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +3.03%  '

# Row 14
$c = $ws.Range('D14')
$c.NumberFormat = '@'
$c.Value = '88.09'
$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  -0.79%  '

# Row 15
$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '0.6670'
$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +0.55%  '

# Row 16
$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '30.539.47'
$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  -0.19%  '

# Row 17
$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +3.56%  '

# Row 18
$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '0.000007846'
$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +0.02%  '

# Row 19
$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -0.02%  '

# Row 20
$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '2.148.21'
$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  +1.14%  '

# Row 21
$c = $ws.Range('B21')
$c.NumberFormat = '@'
$c.Value = 'Uniswap'
$c = $ws.Range('C21')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '5.299'
$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +12.09%  '

# Row 22
$c = $ws.Range('B22')
$c.NumberFormat = '@'
$c.Value = 'BinanceUSD'
$c = $ws.Range('C22')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '1.002'
$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  -0.08%  '

# Row 23
$c = $ws.Range('D23')
$c.NumberFormat = '@'
$c.Value = '189.44'
$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +0.95%  '

# Row 24
$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '6.184'
$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +2.47%  '

# Row 25
$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '9.479'
$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  +2.39%  '

# Row 26
$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '161.73'
$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +2.87%  '

# Row 27
$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '18.33'
$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  +0.48%  '

# Row 28
$c = $ws.Range('D28')
$c.NumberFormat = '@'
$c.Value = '1.930'
$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +5.62%  '

# Row 29
$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '1.474'
$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +5.05%  '

# Row 30
$c = $ws.Range('D30')
$c.NumberFormat = '@'
$c.Value = '4.380'
$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  +3.18%  '

# Row 31
$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '0.09152'
$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  +1.41%  '

# Row 32
$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +4.58%  '

# Row 33
$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '0.05214'
$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +0.75%  '

# Row 34
$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '0.7397'
$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +1.30%  '

# Row 35
$c = $ws.Range('D35')
$c.NumberFormat = '@'
$c.Value = '1.100'
$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +1.91%  '

# Row 36
$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.89%  '

# Row 37
$c = $ws.Range('D37')
$c.NumberFormat = '@'
$c.Value = '0.01832'
$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +0.57%  '

# Row 38
$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '2.680'
$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +0.75%  '

# Row 39
$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '0.9199'
$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +0.20%  '

# Row 40
$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '2.053'
$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +0.45%  '

# Row 41
$c = $ws.Range('D41')
$c.NumberFormat = '@'
$c.Value = '0.4404'
$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  +1.94%  '

# Row 42
$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.932'
$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  +3.72%  '

# Row 43
$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '106.24'
$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +2.36%  '

# Row 44
$c = $ws.Range('D44')
$c.NumberFormat = '@'
$c.Value = '0.9940'
$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  -0.47%  '

# Row 45
$c = $ws.Range('B45')
$c.NumberFormat = '@'
$c.Value = 'Algorand'
$c = $ws.Range('C45')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$c = $ws.Range('D45')
$c.NumberFormat = '@'
$c.Value = '0.1386'
$c = $ws.Range('E45')
$c.NumberFormat = '@'
$c.Value = '  +3.44%  '

# Row 46
$c = $ws.Range('B46')
$c.NumberFormat = '@'
$c.Value = 'Aave'
$c = $ws.Range('C46')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$c = $ws.Range('D46')
$c.NumberFormat = '@'
$c.Value = '68.82'
$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +20.80%  '

# Row 47
$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '7.600'
$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +4.70%  '

# Row 48
$c = $ws.Range('B48')
$c.NumberFormat = '@'
$c.Value = 'Elrond'
$c = $ws.Range('C48')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '35.01'
$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +5.66%  '

# Row 49
$c = $ws.Range('B49')
$c.NumberFormat = '@'
$c.Value = 'EnergySwap'
$c = $ws.Range('C49')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D49')
$c.NumberFormat = '@'
$c.Value = '8.960'
$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +3.88%  '

# Row 50
$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '0.05833'
$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +0.10%  '

# Row 51
$c = $ws.Range('B51')
$c.NumberFormat = '@'
$c.Value = 'Decentraland'
$c = $ws.Range('C51')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '0.3940'
$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -1.84%  '
